# Applies the "Converted figure 3B to greyscale" commit's text changes to the
# Results section discussing the Diversity Outbred (DO) mouse BHB tolerance
# tests.

$d = $word.ActiveDocument

function FindUnique([object]$range, [string]$text) {
    # Executes a forward, case-sensitive, whole-text Find on the supplied
    # range (in place) and returns it collapsed to the found hit.
    $range.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $range
}

# ---------------------------------------------------------------------------
# Change 1: "Diversity outbred mice" -> "Diversity outbred (DO) mice"
# ---------------------------------------------------------------------------

$anchor1 = $d.Content
FindUnique $anchor1 "we performed BHB tolerance tests on diversity outbred mice before or after four weeks of a ketogenic diet.  Diversity outbred mice are genetically unique, so represent the integrated" | Out-Null

$ins1 = $d.Range($anchor1.Start, $anchor1.End)
FindUnique $ins1 "Diversity outbred " | Out-Null
$ins1.Collapse(0)
$ins1.InsertAfter("(DO) ")

Write-Output "Change 1 done"

# ---------------------------------------------------------------------------
# Change 2: rewrite the "within-mouse effects of diet" paragraph
# ---------------------------------------------------------------------------

# 2a: " effects of diet," -> " effects of the diet,"
$r2a = $d.Content
FindUnique $r2a "effects of diet, again showing substantial between-strain variability" | Out-Null
$r2a2 = $d.Range($r2a.Start, $r2a.End)
FindUnique $r2a2 "diet, again" | Out-Null
$r2a2.Collapse(1)
$r2a2.InsertBefore("the ")

Write-Output "Change 2a done"

# 2b: insert "There was more variability in the area under the curve
#     post-diet than pre-diet (<bookmark>).  " right before "Consistent"
$r2b = $d.Content
FindUnique $r2b "genetic differences.  Consistent with our findings from A/J mice" | Out-Null
$r2b2 = $d.Range($r2b.Start, $r2b.End)
FindUnique $r2b2 "Consistent" | Out-Null
$r2b2.Collapse(1)
$r2b2.InsertBefore("There was more variability in the area under the curve post-diet than pre-diet ().  ")

Write-Output "Change 2b done"

# 2c: move the _GoBack bookmark into the new "pre-diet ()" parenthetical
$r2c = $d.Content
FindUnique $r2c "post-diet than pre-diet ().  Consistent" | Out-Null
$r2c2 = $d.Range($r2c.Start, $r2c.End)
FindUnique $r2c2 "pre-diet (" | Out-Null
$r2c2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2c2)

Write-Output "Change 2c done"

# 2d: "with our findings from A/J mice" -> "with our findings using inbred A/J mice"
$r2d = $d.Content
FindUnique $r2d "with our findings from A/J mice the majority of DO mice" | Out-Null
$r2d2 = $d.Range($r2d.Start, $r2d.End)
FindUnique $r2d2 "from" | Out-Null
$r2d2.Text = "using inbred"

Write-Output "Change 2d done"

# 2e: "(XX)" -> "(35 mice)" placeholder fill-in
$r2e = $d.Content
FindUnique $r2e "worsened ketone disposal after diet (XX), with only a small number" | Out-Null
$r2e2 = $d.Range($r2e.Start, $r2e.End)
FindUnique $r2e2 "XX" | Out-Null
$r2e2.Text = "35 mice"

Write-Output "Change 2e done"

# 2f: "(YY)" -> "(10)" placeholder fill-in
$r2f = $d.Content
FindUnique $r2f "improved ketone disposal by our assay (YY)" | Out-Null
$r2f2 = $d.Range($r2f.Start, $r2f.End)
FindUnique $r2f2 "YY" | Out-Null
$r2f2.Text = "10"

Write-Output "Change 2f done"

# 2g: replace the trailing double space at the end of the paragraph with a
#     single space followed by the new "Taking the population..." sentence,
#     including a superscripted "-5" exponent.
$r2g = $d.Content
FindUnique $r2g "by our assay (10).  " | Out-Null
$r2g2 = $d.Range($r2g.Start, $r2g.End)
FindUnique $r2g2 "  " | Out-Null
$r2g2.Text = " "
$r2g2.Collapse(0)
$r2g2.InsertAfter("Taking the population together there was a significant decrease in ketone disposal (increase in baseline adjusted KTT) in these mice (p=4.8 x 10-5 from a paired Wilcoxon test).")

Write-Output "Change 2g done"

# 2h: superscript the "-5" exponent in "p=4.8 x 10-5 from a paired..."
$r2h = $d.Content
FindUnique $r2h "(p=4.8 x 10-5 from a paired Wilcoxon test)." | Out-Null
$r2h2 = $d.Range($r2h.Start, $r2h.End)
FindUnique $r2h2 "-5" | Out-Null
$r2h2.Font.Superscript = $true

Write-Output "Change 2h done"
